$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(7).Insert()
$ws.Columns.Item(7).ColumnWidth = 13.5
